$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.15%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'43.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'6.87%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.832"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'3.29%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08309"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.81%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.780"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.40%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.500"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.95%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.960"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-3.87%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.922"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.05%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9349"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.83%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1253"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.29%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1954"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.19%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09495"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.36%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03976"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'7.17%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'0.81%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001304"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.46%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.005918"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-3.70%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.519"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.56%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'9.044"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'9.50%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1371"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.68%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2572"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-3.06%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04399"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.69%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001255"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.91%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004401"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.51%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'0.69%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003992"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.01%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02799"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.69%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05587"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.61%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007908"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.31%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'0.51%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.009079"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.79%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002101"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-0.63%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'-13.45%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007216"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'4.83%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.15%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003962"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'12.16%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002280"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.18%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.15%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.15%"
$ws.Range("E51").Style = "Normal"
